$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range("D2")
$c.NumberFormat = '@'
$c.Value = '93.079.94'
$c.Style = 'Normal'
$ws.Range("E2").Value = '  +4.88%  '

$c = $ws.Range("D3")
$c.NumberFormat = '@'
$c.Value = '3.283.16'
$c.Style = 'Normal'
$ws.Range("E3").Value = '  +0.07%  '

$ws.Range("E4").Value = '  -0.09%  '

$c = $ws.Range("D5")
$c.NumberFormat = '@'
$c.Value = '219.42'
$c.Style = 'Normal'
$ws.Range("E5").Value = '  +2.93%  '

$c = $ws.Range("D6")
$c.NumberFormat = '@'
$c.Value = '629.91'
$c.Style = 'Normal'
$ws.Range("E6").Value = '  -0.10%  '

$c = $ws.Range("D7")
$c.NumberFormat = '@'
$c.Value = '0.405'
$c.Style = 'Normal'
$ws.Range("E7").Value = '  +3.36%  '

$c = $ws.Range("D8")
$c.NumberFormat = '@'
$c.Value = '0.710'
$c.Style = 'Normal'
$ws.Range("E8").Value = '  +1.45%  '

$ws.Range("E9").Value = '  +0.01%  '

$c = $ws.Range("D10")
$c.NumberFormat = '@'
$c.Value = '3.277.81'
$c.Style = 'Normal'
$ws.Range("E10").Value = '  +0.01%  '

$c = $ws.Range("D11")
$c.NumberFormat = '@'
$c.Value = '0.591'
$c.Style = 'Normal'
$ws.Range("E11").Value = '  +2.12%  '

$c = $ws.Range("D12")
$c.NumberFormat = '@'
$c.Value = '0.0000272'
$c.Style = 'Normal'
$ws.Range("E12").Value = '  +3.11%  '

$ws.Range("E13").Value = '  -3.35%  '

$c = $ws.Range("D14")
$c.NumberFormat = '@'
$c.Value = '34.37'
$c.Style = 'Normal'
$ws.Range("E14").Value = '  +0.44%  '

$c = $ws.Range("D15")
$c.NumberFormat = '@'
$c.Value = '92.560.15'
$c.Style = 'Normal'
$ws.Range("E15").Value = '  +4.77%  '

$c = $ws.Range("D16")
$c.NumberFormat = '@'
$c.Value = '3.862.93'
$c.Style = 'Normal'
$ws.Range("E16").Value = '  -0.40%  '

$c = $ws.Range("D17")
$c.NumberFormat = '@'
$c.Value = '5.34'
$c.Style = 'Normal'
$ws.Range("E17").Value = '  -0.33%  '

$c = $ws.Range("D18")
$c.NumberFormat = '@'
$c.Value = '3.270.04'
$c.Style = 'Normal'
$ws.Range("E18").Value = '  -0.20%  '

$c = $ws.Range("D19")
$c.NumberFormat = '@'
$c.Value = '3.32'
$c.Style = 'Normal'
$ws.Range("E19").Value = '  +6.25%  '

$ws.Range("E20").Value = '  +58.53%  '

$c = $ws.Range("D21")
$c.NumberFormat = '@'
$c.Value = '13.99'
$c.Style = 'Normal'
$ws.Range("E21").Value = '  -0.97%  '

$c = $ws.Range("D22")
$c.NumberFormat = '@'
$c.Value = '448.81'
$c.Style = 'Normal'
$ws.Range("E22").Value = '  +2.83%  '

$c = $ws.Range("D23")
$c.NumberFormat = '@'
$c.Value = '8.85'
$c.Style = 'Normal'
$ws.Range("E23").Value = '  -0.63%  '

$c = $ws.Range("D24")
$c.NumberFormat = '@'
$c.Value = '5.26'
$c.Style = 'Normal'
$ws.Range("E24").Value = '  -2.66%  '

$c = $ws.Range("D25")
$c.NumberFormat = '@'
$c.Value = '5.34'
$c.Style = 'Normal'
$ws.Range("E25").Value = '  +2.78%  '

$c = $ws.Range("D26")
$c.NumberFormat = '@'
$c.Value = '12.15'
$c.Style = 'Normal'
$ws.Range("E26").Value = '  -1.33%  '

$c = $ws.Range("D27")
$c.NumberFormat = '@'
$c.Value = '3.450.03'
$c.Style = 'Normal'
$ws.Range("E27").Value = '  +0.39%  '

$c = $ws.Range("D28")
$c.NumberFormat = '@'
$c.Value = '77.83'
$c.Style = 'Normal'
$ws.Range("E28").Value = '  +0.92%  '

$ws.Range("E29").Value = '  +0.00%  '

$ws.Range("E30").Value = '  -6.27%  '

$c = $ws.Range("D31")
$c.NumberFormat = '@'
$c.Value = '1.00'
$c.Style = 'Normal'
$ws.Range("E31").Value = '  -0.24%  '

$c = $ws.Range("D32")
$c.NumberFormat = '@'
$c.Value = '8.74'
$c.Style = 'Normal'
$ws.Range("E32").Value = '  -1.75%  '

$c = $ws.Range("D33")
$c.NumberFormat = '@'
$c.Value = '554.23'
$c.Style = 'Normal'
$ws.Range("E33").Value = '  -2.98%  '

$c = $ws.Range("D34")
$c.NumberFormat = '@'
$c.Value = '3.85'
$c.Style = 'Normal'
$ws.Range("E34").Value = '  +29.01%  '

$c = $ws.Range("D35")
$c.NumberFormat = '@'
$c.Value = '7.08'
$c.Style = 'Normal'
$ws.Range("E35").Value = '  -0.16%  '

$ws.Range("E36").Value = '  -1.75%  '

$c = $ws.Range("D37")
$c.NumberFormat = '@'
$c.Value = '1.29'
$c.Style = 'Normal'
$ws.Range("E37").Value = '  -8.24%  '

$c = $ws.Range("D38")
$c.NumberFormat = '@'
$c.Value = '22.68'
$c.Style = 'Normal'
$ws.Range("E38").Value = '  +0.05%  '

$c = $ws.Range("D39")
$c.NumberFormat = '@'
$c.Value = '22.50'
$c.Style = 'Normal'
$ws.Range("E39").Value = '  +3.25%  '

$ws.Range("E40").Value = '  -6.44%  '

$ws.Range("E41").Value = '  -0.09%  '

$c = $ws.Range("D42")
$c.NumberFormat = '@'
$c.Value = '0.393'
$c.Style = 'Normal'
$ws.Range("E42").Value = '  -1.48%  '

$c = $ws.Range("D43")
$c.NumberFormat = '@'
$c.Value = '1.99'
$c.Style = 'Normal'
$ws.Range("E43").Value = '  -1.67%  '

$ws.Range("E44").Value = '  -0.07%  '

$c = $ws.Range("D45")
$c.NumberFormat = '@'
$c.Value = '149.85'
$c.Style = 'Normal'
$ws.Range("E45").Value = '  -2.89%  '

$c = $ws.Range("D46")
$c.NumberFormat = '@'
$c.Value = '45.53'
$c.Style = 'Normal'
$ws.Range("E46").Value = '  +1.18%  '

$c = $ws.Range("D47")
$c.NumberFormat = '@'
$c.Value = '178.78'
$c.Style = 'Normal'
$ws.Range("E47").Value = '  -0.95%  '

$ws.Range("E48").Value = '  +1.55%  '

$ws.Range("E49").Value = '  -1.13%  '

$c = $ws.Range("D50")
$c.NumberFormat = '@'
$c.Value = '0.640'
$c.Style = 'Normal'
$ws.Range("E50").Value = '  +2.17%  '

$c = $ws.Range("D51")
$c.NumberFormat = '@'
$c.Value = '4.22'
$c.Style = 'Normal'
$ws.Range("E51").Value = '  -0.47%  '
